# This workbook lists NBA players; the edit reorders the data rows
# (rows 2-18) into a new order while keeping the header row (row 1)
# and the exact same set of player/position/team values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired order for rows 2..18 (Player, Position, Team)
$rows = @(
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Malcolm Brogdon", "PG,SG", "Washington Wizards"),
    @("Dennis Schröder", "PG", "Brooklyn Nets"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("RJ Barrett", "SF,PF", "Toronto Raptors"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Stephon Castle", "PG,SG", "San Antonio Spurs"),
    @("Shai Gilgeous-Alexander", "PG", "Oklahoma City Thunder"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
